# "Updated README and Todo"
# Rename the original sheet to "Features", duplicate it to create a new
# "Bugs" sheet, and populate the Bugs sheet with a fresh set of bug
# reports while keeping the same header/title formatting.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet -------------------------------------------------
$features = $wb.Worksheets.Item(1)
$features.Name = "Features"

# --- Give the 'Create a sniper' row the same Neutral styling used elsewhere ---
$features.Range("B11").Style = "Neutral"

# --- Duplicate the Features sheet (keeps column widths/fonts/fills identical) -
$features.Copy([System.Reflection.Missing]::Value, $features)
$bugs = $wb.Worksheets.Item(2)
$bugs.Name = "Bugs"

# --- The Bugs sheet starts as a clone of Features; strip that down -----------
$bugs.Range("A1:C1").UnMerge()

# Clear out all of the copied Features text below the header row, we'll
# replace it with bug reports.
$bugs.Range("A3:C23").ClearContents()

# Rows 3-19 keep a plain, word-wrapped look (this also strips the
# priority-colour fills those cells inherited from the Features copy).
$bugs.Range("A3:C19").Style = "Normal"

# --- New bug reports (entered in this order so new shared-strings line up) ---
$bugs.Range("C3").Value = "Be able to switch weapons mid-reload."
$bugs.Range("B3").Value = "Prevent player from rotating camera too far."
$bugs.Range("B4").Value = "Re-lock cursor to screen after tabbing out and tabbing back in."
$bugs.Range("A3").Value = "On loss, repeatedly pausing/unpausing causes the games to run for ~1 frame"
$bugs.Range("B5").Value = "If switch weapon before finish reloading, and switch back, the weapon model transform is odd."

$bugs.Range("A3:C19").WrapText = $true

$wb.Save()
